# Applies the "add new mora period (2509) + update totals" edit described
# in the commit: "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row right after the current last data row (60), pushing
#    everything below (the signature block at 65/66) down by one row.
$ws.Rows.Item(61).Insert()

# 2) The newly inserted row 61 should become the new "closing" row of the
#    table (same look the old row 60 had), while row 60 becomes a normal
#    interior row (same look rows 16-59 have). Copy formats accordingly.
$ws.Rows.Item(60).Copy()
$ws.Rows.Item(61).PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(59).Copy()
$ws.Rows.Item(60).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Fill in the data for the new period (2509) - same worker, same amounts
#    as the rest of the table.
$ws.Range("B61").Value = "CC"
$ws.Range("C61").Value = "9176732"
$ws.Range("D61").Value = "OSCAR ENRIQUE ALVIS MELENDEZ"
$ws.Range("E61").Value = "2509"
$ws.Range("F61").Value = 31320
$ws.Range("G61").Value = 783000
$ws.Range("H61").Value = ""
$ws.Range("I61").Value = ""
$ws.Range("J61").Value = ""

# 4) Update the summary fields: one more mora period, and the total mora
#    value grows by the new period's "Valor Mora" (31320).
$ws.Range("E11").Value = 1449400
$ws.Range("F13").Value = 45
